$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# Set updated values
$ws1.Range("H2").Value = 589.35297
$ws1.Range("I2").Value = 652.5714
$ws1.Range("J2").Value = 294.33334
$ws1.Range("K2").Value = 652.5714
$ws1.Range("L2").Value = 294.33334
$ws1.Range("M2").Value = -539.5714
$ws1.Range("N2").Value = -520.33334
$ws1.Range("H6").Value = 37.333332
$ws1.Range("I6").Value = 31
$ws1.Range("J6").Value = 50
$ws1.Range("K6").Value = 93
$ws1.Range("L6").Value = 150
$ws1.Range("M6").Value = 19
$ws1.Range("N6").Value = -374
$ws1.Range("H12").Value = 724.55
$ws1.Range("I12").Value = 270.17648
$ws1.Range("J12").Value = 3299.3333
$ws1.Range("K12").Value = 270.17648
$ws1.Range("L12").Value = 3299.3333
$ws1.Range("M12").Value = -100.17648
$ws1.Range("N12").Value = -3639.3333
$ws1.Range("H15").Value = 132275.4
$ws1.Range("I15").Value = 132275.4
$ws1.Range("K15").Value = 396826.2
$ws1.Range("M15").Value = -396657.2
$ws1.Range("H17").Value = 1871.742
$ws1.Range("I17").Value = 1132.8334
$ws1.Range("J17").Value = 2049.08
$ws1.Range("K17").Value = 3398.5002
$ws1.Range("L17").Value = 6147.24
$ws1.Range("M17").Value = -3230.5002
$ws1.Range("N17").Value = -6483.24
$ws1.Range("H21").Value = 51333.332
$ws1.Range("I21").Value = 51333.332
$ws1.Range("J21").Value = 0
$ws1.Range("K21").Value = 51333.332
$ws1.Range("L21").Value = 0
$ws1.Range("N21").Value = -50865.332
$ws1.Range("H23").Value = 51333.332
$ws1.Range("I23").Value = 51333.332
$ws1.Range("J23").Value = 0
$ws1.Range("K23").Value = 51333.332
$ws1.Range("L23").Value = 0
$ws1.Range("N23").Value = -51099.332
$ws1.Range("H29").Value = 2963.818
$ws1.Range("J29").Value = 4074.125
$ws1.Range("L29").Value = 12222.375
$ws1.Range("N29").Value = -12784.375
$ws1.Range("H33").Value = 646.3077
$ws1.Range("I33").Value = 413.25
$ws1.Range("J33").Value = 1019.2
$ws1.Range("K33").Value = 413.25
$ws1.Range("L33").Value = 1019.2
$ws1.Range("M33").Value = -184.25
$ws1.Range("N33").Value = -1477.2
$ws1.Range("H39").Value = 36.833332
$ws1.Range("I39").Value = 42.4
$ws1.Range("J39").Value = 9
$ws1.Range("K39").Value = 127.2
$ws1.Range("L39").Value = 27
$ws1.Range("M39").Value = 168.8
$ws1.Range("N39").Value = -619
$ws1.Range("H41").Value = 271
$ws1.Range("I41").Value = 238
$ws1.Range("J41").Value = 274.3
$ws1.Range("K41").Value = 238
$ws1.Range("L41").Value = 274.3
$ws1.Range("M41").Value = 202
$ws1.Range("N41").Value = -1154.3
$ws1.Range("H43").Value = 1690.6666
$ws1.Range("I43").Value = 883.1667
$ws1.Range("J43").Value = 2498.1667
$ws1.Range("K43").Value = 883.1667
$ws1.Range("L43").Value = 2498.1667
$ws1.Range("M43").Value = -814.1667
$ws1.Range("N43").Value = -2636.1667
$ws1.Range("H55").Value = 374.5
$ws1.Range("I55").Value = 147.5
$ws1.Range("J55").Value = 525.8333
$ws1.Range("K55").Value = 147.5
$ws1.Range("L55").Value = 525.8333
$ws1.Range("M55").Value = 66.5
$ws1.Range("N55").Value = -953.8333
$ws1.Range("H61").Value = 715
$ws1.Range("I61").Value = 715
$ws1.Range("K61").Value = 2145
$ws1.Range("M61").Value = -1973
$ws1.Range("H62").Value = 11911919
$ws1.Range("I62").Value = 20842330
$ws1.Range("K62").Value = 20842330
$ws1.Range("M62").Value = -20841706
$ws1.Range("H65").Value = 11911919
$ws1.Range("I65").Value = 20842330
$ws1.Range("K65").Value = 104211650
$ws1.Range("M65").Value = -104208530
$ws1.Range("H100").Value = 5666
$ws1.Range("I100").Value = 5207.5
$ws1.Range("K100").Value = 5207.5
$ws1.Range("M100").Value = -4666.5
$ws1.Range("H132").Value = 226634.62
$ws1.Range("I132").Value = 259464.67
$ws1.Range("J132").Value = 9135.5
$ws1.Range("K132").Value = 778394.01
$ws1.Range("L132").Value = 27406.5
$ws1.Range("M132").Value = -775864.01
$ws1.Range("N132").Value = -32466.5
$ws1.Range("H137").Value = 5226.5264
$ws1.Range("I137").Value = 5255
$ws1.Range("K137").Value = 15765
$ws1.Range("M137").Value = -13215
$ws1.Range("H138").Value = 3663.311
$ws1.Range("I138").Value = 2014.7368
$ws1.Range("J138").Value = 4868.0386
$ws1.Range("K138").Value = 6044.2104
$ws1.Range("L138").Value = 14604.1158
$ws1.Range("M138").Value = -904.2103999999999
$ws1.Range("N138").Value = -24884.1158
$ws1.Range("H141").Value = 2247.5
$ws1.Range("I141").Value = 1854.2858
$ws1.Range("J141").Value = 5000
$ws1.Range("K141").Value = 5562.857400000001
$ws1.Range("L141").Value = 15000
$ws1.Range("M141").Value = -382.8574000000008
$ws1.Range("N141").Value = -25360
$ws2.Range("H2").Value = 4386.2144
$ws2.Range("I2").Value = 1886.7222
$ws2.Range("K2").Value = 1886.7222
$ws2.Range("M2").Value = -1773.7222
$ws2.Range("H5").Value = 239.6
$ws2.Range("I5").Value = 199.61539
$ws2.Range("K5").Value = 199.61539
$ws2.Range("M5").Value = -87.61538999999999
$ws2.Range("H32").Value = 55569556
$ws2.Range("I32").Value = 0
$ws2.Range("K32").Value = 0
$ws2.Range("H45").Value = 1655.3334
$ws2.Range("I45").Value = 916.75
$ws2.Range("K45").Value = 916.75
$ws2.Range("M45").Value = -539.75
$ws2.Range("H55").Value = 0
$ws2.Range("J55").Value = 0
$ws2.Range("N55").Value = 0
$ws2.Range("H61").Value = 4921.6562
$ws2.Range("I61").Value = 5643.579
$ws2.Range("J61").Value = 4616.844
$ws2.Range("K61").Value = 5643.579
$ws2.Range("L61").Value = 4616.844
$ws2.Range("M61").Value = -5431.579
$ws2.Range("N61").Value = -5040.844
$ws2.Range("H74").Value = 4334.0312
$ws2.Range("I74").Value = 3796.087
$ws2.Range("J74").Value = 5708.778
$ws2.Range("K74").Value = 3796.087
$ws2.Range("L74").Value = 5708.778
$ws2.Range("M74").Value = -2922.087
$ws2.Range("N74").Value = -7456.778
$ws2.Range("H77").Value = 4334.0312
$ws2.Range("I77").Value = 3796.087
$ws2.Range("J77").Value = 5708.778
$ws2.Range("K77").Value = 18980.435
$ws2.Range("L77").Value = 28543.89
$ws2.Range("M77").Value = -14612.435
$ws2.Range("N77").Value = -37279.89
$ws2.Range("H80").Value = 11172.5
$ws2.Range("I80").Value = 10000
$ws2.Range("K80").Value = 10000
$ws2.Range("M80").Value = -9002
$ws2.Range("H83").Value = 11172.5
$ws2.Range("I83").Value = 10000
$ws2.Range("K83").Value = 30000
$ws2.Range("M83").Value = -25008
$ws2.Range("H102").Value = 2358.75
$ws2.Range("I102").Value = 2195.8572
$ws2.Range("J102").Value = 3499
$ws2.Range("K102").Value = 2195.8572
$ws2.Range("L102").Value = 3499
$ws2.Range("M102").Value = -573.8571999999999
$ws2.Range("N102").Value = -6743
$ws2.Range("H110").Value = 5435.5
$ws2.Range("I110").Value = 2350.6
$ws2.Range("K110").Value = 2350.6
$ws2.Range("M110").Value = -305.5999999999999
$ws2.Range("H116").Value = 4386.2144
$ws2.Range("I116").Value = 1886.7222
$ws2.Range("K116").Value = 1886.7222
$ws2.Range("M116").Value = 407.2778000000001
$ws2.Range("H122").Value = 3332.32
$ws2.Range("I122").Value = 2700.3635
$ws2.Range("K122").Value = 8101.0905
$ws2.Range("M122").Value = -5651.0905
$ws2.Range("H132").Value = 589017.9399999999
$ws2.Range("I132").Value = 715829.5600000001
$ws2.Range("J132").Value = 81771.38
$ws2.Range("K132").Value = 2147488.68
$ws2.Range("L132").Value = 245314.14
$ws2.Range("M132").Value = -2144958.68
$ws2.Range("N132").Value = -250374.14
$ws2.Range("H136").Value = 4921.6562
$ws2.Range("I136").Value = 5643.579
$ws2.Range("J136").Value = 4616.844
$ws2.Range("K136").Value = 16930.737
$ws2.Range("L136").Value = 13850.532
$ws2.Range("M136").Value = -14380.737
$ws2.Range("N136").Value = -18950.532
$ws3.Range("H3").Value = 4386.2144
$ws3.Range("I3").Value = 1886.7222
$ws3.Range("K3").Value = 1886.7222
$ws3.Range("M3").Value = -1772.7222
$ws3.Range("H4").Value = 239.6
$ws3.Range("I4").Value = 199.61539
$ws3.Range("K4").Value = 199.61539
$ws3.Range("M4").Value = -84.61538999999999
$ws3.Range("H94").Value = 3934.6843
$ws3.Range("I94").Value = 2443.15
$ws3.Range("K94").Value = 2443.15
$ws3.Range("M94").Value = -1992.15
$ws3.Range("H105").Value = 1983.8235
$ws3.Range("I105").Value = 2020.3125
$ws3.Range("J105").Value = 1400
$ws3.Range("K105").Value = 2020.3125
$ws3.Range("L105").Value = 1400
$ws3.Range("M105").Value = -273.3125
$ws3.Range("N105").Value = -4894
$ws3.Range("H107").Value = 6260632
$ws3.Range("I107").Value = 7700777.5
$ws3.Range("K107").Value = 7700777.5
$ws3.Range("M107").Value = -7698857.5
$ws3.Range("H134").Value = 1124984.2
$ws3.Range("I134").Value = 1466217.1
$ws3.Range("J134").Value = 8222.272000000001
$ws3.Range("K134").Value = 4398651.300000001
$ws3.Range("L134").Value = 24666.816
$ws3.Range("M134").Value = -4396116.300000001
$ws3.Range("N134").Value = -29736.816
$ws4.Range("H16").Value = 35721384
$ws4.Range("I16").Value = 55561260
$ws4.Range("J16").Value = 9602.4
$ws4.Range("K16").Value = 55561260
$ws4.Range("L16").Value = 9602.4
$ws4.Range("M16").Value = -55560973
$ws4.Range("N16").Value = -10176.4
$ws4.Range("H31").Value = 18530260
$ws4.Range("I31").Value = 45473276
$ws4.Range("K31").Value = 45473276
$ws4.Range("M31").Value = -45472981
$ws4.Range("H34").Value = 18530260
$ws4.Range("I34").Value = 45473276
$ws4.Range("K34").Value = 45473276
$ws4.Range("M34").Value = -45473074
$ws4.Range("H58").Value = 71448200
$ws4.Range("I58").Value = 111127256
$ws4.Range("J58").Value = 25899.6
$ws4.Range("K58").Value = 111127256
$ws4.Range("L58").Value = 25899.6
$ws4.Range("M58").Value = -111127053
$ws4.Range("N58").Value = -26305.6
$ws4.Range("H113").Value = 35721384
$ws4.Range("I113").Value = 55561260
$ws4.Range("J113").Value = 9602.4
$ws4.Range("K113").Value = 55561260
$ws4.Range("L113").Value = 9602.4
$ws4.Range("M113").Value = -55559090
$ws4.Range("N113").Value = -13942.4
$ws4.Range("H132").Value = 8234
$ws4.Range("I132").Value = 4489.0415
$ws4.Range("J132").Value = 21073.857
$ws4.Range("K132").Value = 13467.1245
$ws4.Range("L132").Value = 63221.571
$ws4.Range("M132").Value = -10937.1245
$ws4.Range("N132").Value = -68281.571
$ws4.Range("H134").Value = 55565976
$ws4.Range("I134").Value = 76928060
$ws4.Range("J134").Value = 24549.6
$ws4.Range("K134").Value = 230784180
$ws4.Range("L134").Value = 73648.79999999999
$ws4.Range("M134").Value = -230781645
$ws4.Range("N134").Value = -78718.79999999999
$ws4.Range("H136").Value = 71448200
$ws4.Range("I136").Value = 111127256
$ws4.Range("J136").Value = 25899.6
$ws4.Range("K136").Value = 333381768
$ws4.Range("L136").Value = 77698.79999999999
$ws4.Range("M136").Value = -333379218
$ws4.Range("N136").Value = -82798.79999999999
$ws4.Range("H138").Value = 87827.60000000001
$ws4.Range("J138").Value = 79712.664
$ws4.Range("L138").Value = 79712.664
$ws4.Range("N138").Value = -89992.664
$ws4.Range("H140").Value = 40780
$ws4.Range("J140").Value = 40780
$ws4.Range("L140").Value = 40780
$ws4.Range("N140").Value = -51140
$ws5.Range("H2").Value = 87.56521600000001
$ws5.Range("I2").Value = 37.909092
$ws5.Range("J2").Value = 133.08333
$ws5.Range("K2").Value = 227.454552
$ws5.Range("L2").Value = 798.4999799999999
$ws5.Range("M2").Value = -114.454552
$ws5.Range("N2").Value = -1024.49998
$ws5.Range("H5").Value = 952.3333
$ws5.Range("I5").Value = 563.8333
$ws5.Range("J5").Value = 1729.3334
$ws5.Range("K5").Value = 1691.4999
$ws5.Range("L5").Value = 5188.0002
$ws5.Range("M5").Value = -1579.4999
$ws5.Range("N5").Value = -5412.0002
$ws5.Range("H23").Value = 79.5
$ws5.Range("I23").Value = 53.2
$ws5.Range("J23").Value = 98.28570999999999
$ws5.Range("K23").Value = 159.6
$ws5.Range("L23").Value = 294.85713
$ws5.Range("M23").Value = 75.39999999999998
$ws5.Range("N23").Value = -764.85713
$ws5.Range("H97").Value = 1109
$ws5.Range("J97").Value = 1385.1428
$ws5.Range("L97").Value = 4155.428400000001
$ws5.Range("N97").Value = -5147.428400000001
$ws5.Range("H98").Value = 1084.6154
$ws5.Range("J98").Value = 1139.5
$ws5.Range("L98").Value = 3418.5
$ws5.Range("N98").Value = -6414.5
$ws5.Range("H107").Value = 3356.1904
$ws5.Range("I107").Value = 686.5
$ws5.Range("J107").Value = 4190.4688
$ws5.Range("K107").Value = 2059.5
$ws5.Range("L107").Value = 12571.4064
$ws5.Range("M107").Value = -139.5
$ws5.Range("N107").Value = -16411.4064
$ws5.Range("H114").Value = 1320.1
$ws5.Range("I114").Value = 584.5
$ws5.Range("K114").Value = 1753.5
$ws5.Range("M114").Value = 1500.5
$ws5.Range("H117").Value = 2355.7144
$ws5.Range("I117").Value = 869.875
$ws5.Range("K117").Value = 2609.625
$ws5.Range("M117").Value = 832.375
$ws5.Range("H122").Value = 3490.6206
$ws5.Range("I122").Value = 987
$ws5.Range("K122").Value = 8883
$ws5.Range("M122").Value = -6433
$ws5.Range("H135").Value = 952.3333
$ws5.Range("I135").Value = 563.8333
$ws5.Range("J135").Value = 1729.3334
$ws5.Range("K135").Value = 5074.4997
$ws5.Range("L135").Value = 15564.0006
$ws5.Range("M135").Value = -2539.4997
$ws5.Range("N135").Value = -20634.0006
$ws5.Range("H137").Value = 3474.0715
$ws5.Range("I137").Value = 2495.75
$ws5.Range("J137").Value = 3865.4
$ws5.Range("K137").Value = 7487.25
$ws5.Range("L137").Value = 11596.2
$ws5.Range("M137").Value = -2387.25
$ws5.Range("N137").Value = -21796.2
$ws6.Range("H2").Value = 72.25
$ws6.Range("I2").Value = 54.5
$ws6.Range("J2").Value = 90
$ws6.Range("K2").Value = 54.5
$ws6.Range("L2").Value = 90
$ws6.Range("M2").Value = 58.5
$ws6.Range("N2").Value = -316
$ws6.Range("H80").Value = 7388
$ws6.Range("I80").Value = 5989.3335
$ws6.Range("J80").Value = 7912.5
$ws6.Range("K80").Value = 5989.3335
$ws6.Range("L80").Value = 7912.5
$ws6.Range("M80").Value = -4991.3335
$ws6.Range("N80").Value = -9908.5
$ws6.Range("H83").Value = 7388
$ws6.Range("I83").Value = 5989.3335
$ws6.Range("J83").Value = 7912.5
$ws6.Range("K83").Value = 29946.6675
$ws6.Range("L83").Value = 39562.5
$ws6.Range("M83").Value = -24954.6675
$ws6.Range("N83").Value = -49546.5
$ws6.Range("H95").Value = 21400
$ws6.Range("J95").Value = 21400
$ws6.Range("L95").Value = 21400
$ws6.Range("N95").Value = -26892
$ws6.Range("H97").Value = 2457.2
$ws6.Range("I97").Value = 2129.7334
$ws6.Range("J97").Value = 3439.6
$ws6.Range("K97").Value = 2129.7334
$ws6.Range("L97").Value = 3439.6
$ws6.Range("M97").Value = -1633.7334
$ws6.Range("N97").Value = -4431.6
$ws6.Range("H99").Value = 13085.5
$ws6.Range("I99").Value = 13085.5
$ws6.Range("J99").Value = 0
$ws6.Range("K99").Value = 13085.5
$ws6.Range("L99").Value = 0
$ws6.Range("N99").Value = -10839.5
$ws6.Range("H102").Value = 1103469.1
$ws6.Range("I102").Value = 2050189
$ws6.Range("K102").Value = 2050189
$ws6.Range("M102").Value = -2048567
$ws6.Range("H107").Value = 787.4
$ws6.Range("J107").Value = 1182.7778
$ws6.Range("L107").Value = 1182.7778
$ws6.Range("N107").Value = -5022.7778
$ws6.Range("H113").Value = 6847.1177
$ws6.Range("I113").Value = 1459.625
$ws6.Range("K113").Value = 1459.625
$ws6.Range("M113").Value = 710.375
$ws6.Range("H126").Value = 38475836
$ws6.Range("I126").Value = 83339656
$ws6.Range("J126").Value = 21128.285
$ws6.Range("K126").Value = 250018968
$ws6.Range("L126").Value = 63384.855
$ws6.Range("M126").Value = -250016498
$ws6.Range("N126").Value = -68324.855
$ws6.Range("H132").Value = 5416.4883
$ws6.Range("I132").Value = 5687.6562
$ws6.Range("J132").Value = 4627.636
$ws6.Range("K132").Value = 17062.9686
$ws6.Range("L132").Value = 13882.908
$ws6.Range("M132").Value = -14532.9686
$ws6.Range("N132").Value = -18942.908
$ws7.Range("H16").Value = 6989.636
$ws7.Range("I16").Value = 3982.3333
$ws7.Range("K16").Value = 3982.3333
$ws7.Range("M16").Value = -3812.3333
$ws7.Range("H22").Value = 1261
$ws7.Range("I22").Value = 1146.3334
$ws7.Range("J22").Value = 1347
$ws7.Range("K22").Value = 1146.3334
$ws7.Range("L22").Value = 1347
$ws7.Range("M22").Value = -851.3334
$ws7.Range("N22").Value = -1937
$ws7.Range("H27").Value = 1261
$ws7.Range("I27").Value = 1146.3334
$ws7.Range("J27").Value = 1347
$ws7.Range("K27").Value = 1146.3334
$ws7.Range("L27").Value = 1347
$ws7.Range("M27").Value = -1039.3334
$ws7.Range("N27").Value = -1561
$ws7.Range("H40").Value = 4795.231
$ws7.Range("I40").Value = 3894.4546
$ws7.Range("J40").Value = 9749.5
$ws7.Range("K40").Value = 3894.4546
$ws7.Range("L40").Value = 9749.5
$ws7.Range("M40").Value = -3758.4546
$ws7.Range("N40").Value = -10021.5
$ws7.Range("H46").Value = 166667550
$ws7.Range("I46").Value = 1300
$ws7.Range("J46").Value = 250000670
$ws7.Range("K46").Value = 1300
$ws7.Range("L46").Value = 250000670
$ws7.Range("M46").Value = -1112
$ws7.Range("N46").Value = -250001046
$ws7.Range("H61").Value = 5096
$ws7.Range("I61").Value = 4424.4546
$ws7.Range("J61").Value = 19870
$ws7.Range("K61").Value = 4424.4546
$ws7.Range("L61").Value = 19870
$ws7.Range("M61").Value = -4222.4546
$ws7.Range("N61").Value = -20274
$ws7.Range("H68").Value = 3633.0557
$ws7.Range("I68").Value = 2426.4
$ws7.Range("J68").Value = 9666.333000000001
$ws7.Range("K68").Value = 2426.4
$ws7.Range("L68").Value = 9666.333000000001
$ws7.Range("M68").Value = -1677.4
$ws7.Range("N68").Value = -11164.333
$ws7.Range("H71").Value = 3633.0557
$ws7.Range("I71").Value = 2426.4
$ws7.Range("J71").Value = 9666.333000000001
$ws7.Range("K71").Value = 12132
$ws7.Range("L71").Value = 48331.665
$ws7.Range("M71").Value = -8388
$ws7.Range("N71").Value = -55819.665
$ws7.Range("H113").Value = 5096
$ws7.Range("I113").Value = 4424.4546
$ws7.Range("J113").Value = 19870
$ws7.Range("K113").Value = 4424.4546
$ws7.Range("L113").Value = 19870
$ws7.Range("M113").Value = -2254.4546
$ws7.Range("N113").Value = -24210
$ws7.Range("H119").Value = 82987.86
$ws7.Range("J119").Value = 82987.86
$ws7.Range("L119").Value = 82987.86
$ws7.Range("N119").Value = -92663.86
$ws7.Range("H122").Value = 1999899.4
$ws7.Range("I122").Value = 3328932.2
$ws7.Range("J122").Value = 6350
$ws7.Range("K122").Value = 9986796.600000001
$ws7.Range("L122").Value = 19050
$ws7.Range("M122").Value = -9984346.600000001
$ws7.Range("N122").Value = -23950
$ws7.Range("H132").Value = 7939586.5
$ws7.Range("I132").Value = 15875215
$ws7.Range("J132").Value = 3957.889
$ws7.Range("K132").Value = 47625645
$ws7.Range("L132").Value = 11873.667
$ws7.Range("M132").Value = -47623115
$ws7.Range("N132").Value = -16933.667
$ws7.Range("H136").Value = 42745444
$ws7.Range("I136").Value = 79374110
$ws7.Range("J136").Value = 11991.167
$ws7.Range("K136").Value = 238122330
$ws7.Range("L136").Value = 35973.501
$ws7.Range("M136").Value = -238119780
$ws7.Range("N136").Value = -41073.501
$ws8.Range("H51").Value = 23899.143
$ws8.Range("J51").Value = 24306
$ws8.Range("L51").Value = 24306
$ws8.Range("N51").Value = -25326
$ws8.Range("H100").Value = 1915.1333
$ws8.Range("J100").Value = 3451.3333
$ws8.Range("L100").Value = 6902.6666
$ws8.Range("N100").Value = -7984.6666
$ws8.Range("H122").Value = 9105.714
$ws8.Range("I122").Value = 3669.889
$ws8.Range("K122").Value = 11009.667
$ws8.Range("M122").Value = -8559.667000000001
$ws8.Range("H123").Value = 62631.668
$ws8.Range("J123").Value = 62631.668
$ws8.Range("L123").Value = 62631.668
$ws8.Range("N123").Value = -72431.66800000001
$ws8.Range("H132").Value = 3929901.8
$ws8.Range("I132").Value = 5336.794
$ws8.Range("J132").Value = 11779032
$ws8.Range("K132").Value = 16010.382
$ws8.Range("L132").Value = 35337096
$ws8.Range("M132").Value = -13480.382
$ws8.Range("N132").Value = -35342156
$ws8.Range("H135").Value = 99444
$ws8.Range("J135").Value = 99444
$ws8.Range("L135").Value = 99444
$ws8.Range("N135").Value = -109584
$ws8.Range("H136").Value = 15163152
$ws8.Range("I136").Value = 20843416
$ws8.Range("K136").Value = 62530248
$ws8.Range("M136").Value = -62527698

# Clear cells that no longer exist after the edit
$ws1.Range("M21").ClearContents()
$ws1.Range("M23").ClearContents()
$ws2.Range("M32").ClearContents()
$ws2.Range("L55").ClearContents()
$ws6.Range("M99").ClearContents()

Write-Output "Applied 534 cell updates and 5 clears"